$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.799.63"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "1.810.90"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'310.38"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.4950"
$ws.Range("E7").Value = "  -5.59%  "
$ws.Range("D8").Value = "'0.3885"
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("D9").Value = "'0.09431"
$ws.Range("E9").Value = "  +19.57%  "
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("D11").Value = "'40.45"
$ws.Range("E11").Value = "  -2.83%  "
$ws.Range("D12").Value = "'6.412"
$ws.Range("E12").Value = "  +2.26%  "
$ws.Range("D13").Value = "'1.003"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").Value = "'20.44"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.806.40"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "'7.276"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "'0.00001124"
$ws.Range("E17").Value = "  +3.80%  "
$ws.Range("D18").Value = "'92.96"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "'0.06618"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "'17.09"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("D22").Value = "'5.942"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "27.869.01"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").Value = "'11.12"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").Value = "'2.252"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("D26").Value = "'20.65"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "'156.67"
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("D28").Value = "2.019.97"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").Value = "'2.423"
$ws.Range("E29").Value = "  +3.61%  "
$ws.Range("D30").Value = "'127.51"
$ws.Range("E30").Value = "  +3.24%  "
$ws.Range("D31").Value = "'0.1064"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "'1.046"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").Value = "'5.555"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").Value = "'3.635"
$ws.Range("E34").Value = "  -1.25%  "
$ws.Range("D35").Value = "'0.06782"
$ws.Range("E35").Value = "  -5.26%  "
$ws.Range("D36").Value = "'9.111"
$ws.Range("E36").Value = "  +5.26%  "
$ws.Range("D37").Value = "'0.02314"
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("D38").Value = "'0.2141"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'4.942"
$ws.Range("E39").Value = "  -2.37%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'11.29"
$ws.Range("E40").Value = "  -5.92%  "
$ws.Range("D41").Value = "'0.6184"
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("D42").Value = "'1.002"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").Value = "'1.142"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.296"
$ws.Range("E45").Value = "  -5.46%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5877"
$ws.Range("E46").Value = "  -5.46%  "
$ws.Range("D47").Value = "'3.703"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("D48").Value = "'123.82"
$ws.Range("E48").Value = "  -2.83%  "
$ws.Range("D49").Value = "'1.944"
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("E50").Value = "  -4.67%  "
$ws.Range("D51").Value = "'0.06779"
$ws.Range("E51").Value = "  +0.12%  "
